$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.532.30'
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").Value = '1.569.92'
$ws.Range("E3").Value = '  -1.66%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.18%  '

$ws.Range("E6").Value = '  -0.74%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.20'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("E10").Value = '  -1.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0591'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0887'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '1.793.98'
$ws.Range("E13").Value = '  -1.71%  '

$ws.Range("D14").Value = '1.575.00'
$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.520'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.31%  '

$ws.Range("D16").Value = '28.497.37'
$ws.Range("E16").Value = '  +0.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.09%  '

$ws.Range("E21").Value = '  -2.91%  '

$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.86'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.64%  '

$ws.Range("E25").Value = '  +7.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("E27").Value = '  -1.32%  '

$ws.Range("E28").Value = '  -2.79%  '

$ws.Range("E29").Value = '  -3.82%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0483'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.38%  '

$ws.Range("E32").Value = '  -3.68%  '

$ws.Range("E33").Value = '  -1.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.29%  '

$ws.Range("D35").Value = '1.393.87'
$ws.Range("E35").Value = '  -0.45%  '

$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.53'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.39%  '

$ws.Range("E38").Value = '  +0.51%  '

$ws.Range("E39").Value = '  +3.37%  '

$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.521'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.99%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("E44").Value = '  -3.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0467'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.971'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '62.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.68%  '

$ws.Range("E49").Value = '  -1.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.52%  '

$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  -1.20%  '
